$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a literal text value into a cell, bypassing Excel's
# automatic number inference for strings that look numeric (e.g. "210.32"),
# then restore the default "Normal" style so no stray formatting is left behind.
function Set-TextValue($addr, $value) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

$ws.Range('D2').Value = '28.301.49'
$ws.Range('E2').Value = '  -0.97%  '
$ws.Range('D3').Value = '1.561.03'
$ws.Range('E3').Value = '  -0.36%  '
Set-TextValue 'D5' '210.32'
$ws.Range('E5').Value = '  -0.67%  '
Set-TextValue 'D6' '0.491'
$ws.Range('E6').Value = '  -0.50%  '
$ws.Range('E7').Value = '  +0.08%  '
Set-TextValue 'D8' '44.33'
$ws.Range('E8').Value = '  -4.48%  '
Set-TextValue 'D9' '23.64'
$ws.Range('E9').Value = '  -2.55%  '
$ws.Range('E10').Value = '  -1.63%  '
$ws.Range('E11').Value = '  -0.95%  '
Set-TextValue 'D12' '0.0893'
$ws.Range('E12').Value = '  +0.91%  '
$ws.Range('D13').Value = '1.784.24'
$ws.Range('E13').Value = '  -0.32%  '
$ws.Range('D14').Value = '1.563.78'
$ws.Range('E14').Value = '  -0.22%  '
$ws.Range('D15').Value = '28.296.85'
$ws.Range('E16').Value = '  -0.79%  '
Set-TextValue 'D17' '0.512'
$ws.Range('E17').Value = '  -1.95%  '
$ws.Range('E18').Value = '  -1.93%  '
Set-TextValue 'D19' '227.88'
$ws.Range('E19').Value = '  -0.41%  '
Set-TextValue 'D20' '7.35'
$ws.Range('E20').Value = '  +0.06%  '
$ws.Range('D21').Value = '0.0₃0678'
$ws.Range('E21').Value = '  -2.42%  '
$ws.Range('E22').Value = '  +0.05%  '
$ws.Range('E23').Value = '  +1.17%  '
Set-TextValue 'D24' '8.88'
$ws.Range('E24').Value = '  -3.04%  '
Set-TextValue 'D25' '2.04'
$ws.Range('E25').Value = '  -2.36%  '
Set-TextValue 'D26' '150.16'
$ws.Range('E26').Value = '  -0.42%  '
$ws.Range('E28').Value = '  -0.55%  '
$ws.Range('E29').Value = '  -2.02%  '
$ws.Range('E30').Value = '  +0.03%  '
Set-TextValue 'D31' '0.0476'
$ws.Range('E31').Value = '  +2.22%  '
$ws.Range('E32').Value = '  -3.04%  '
$ws.Range('E33').Value = '  -1.41%  '
$ws.Range('E34').Value = '  -1.36%  '
$ws.Range('D35').Value = '1.378.06'
$ws.Range('E35').Value = '  -1.48%  '
$ws.Range('E36').Value = '  +1.23%  '
$ws.Range('E37').Value = '  -3.85%  '
$ws.Range('E38').Value = '  -0.35%  '
$ws.Range('E39').Value = '  +2.22%  '
$ws.Range('E40').Value = '  -2.19%  '
$ws.Range('E41').Value = '  -3.31%  '
Set-TextValue 'D42' '1.92'
$ws.Range('E42').Value = '  +2.37%  '
$ws.Range('E43').Value = '  +0.02%  '
Set-TextValue 'D44' '0.0471'
$ws.Range('E44').Value = '  -0.85%  '
$ws.Range('E45').Value = '  -1.41%  '
$ws.Range('E46').Value = '  -3.61%  '
Set-TextValue 'D47' '62.01'
$ws.Range('E47').Value = '  -1.33%  '
Set-TextValue 'D48' '0.918'
$ws.Range('E48').Value = '  -5.92%  '
$ws.Range('D49').Value = '1.697.72'
$ws.Range('E49').Value = '  -0.19%  '
Set-TextValue 'D50' '85.29'
$ws.Range('E50').Value = '  -1.25%  '
$ws.Range('E51').Value = '  -2.27%  '
